# Update NATMI LR-pair edge-weight metrics (Epha4-Efnb1) with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.030956000000001
$ws.Range("H2").Value = 18.092868
$ws.Range("I2").Value = 0.364814105361131
$ws.Range("J2").Value = 0.3648141053611309
$ws.Range("M2").Value = 8.841467
$ws.Range("N2").Value = 26.524401
$ws.Range("O2").Value = 0.5917001192060068
$ws.Range("P2").Value = 0.5917001192060067
$ws.Range("Q2").Value = 53.322498452452
$ws.Range("R2").Value = 479.9024860720681
$ws.Range("S2").Value = 0.2158605496302139
$ws.Range("T2").Value = 0.2158605496302139

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.030956000000001
$ws.Range("H3").Value = 18.092868
$ws.Range("I3").Value = 0.364814105361131
$ws.Range("J3").Value = 0.3648141053611309
$ws.Range("O3").Value = 0.2746155987184545
$ws.Range("P3").Value = 0.2746155987184545
$ws.Range("Q3").Value = 24.747654026728
$ws.Range("R3").Value = 222.728886240552
$ws.Range("S3").Value = 0.1001836439646843
$ws.Range("T3").Value = 0.1001836439646843

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.030956000000001
$ws.Range("H4").Value = 18.092868
$ws.Range("I4").Value = 0.364814105361131
$ws.Range("J4").Value = 0.3648141053611309
$ws.Range("M4").Value = 1.997574666666667
$ws.Range("N4").Value = 5.992724
$ws.Range("O4").Value = 0.1336842820755386
$ws.Range("P4").Value = 0.1336842820755386
$ws.Range("Q4").Value = 12.04728492138134
$ws.Range("R4").Value = 108.425564292432
$ws.Range("S4").Value = 0.04876991176623269
$ws.Range("T4").Value = 0.04876991176623268

# Row 5
$ws.Range("I5").Value = 0.4107214552505144
$ws.Range("J5").Value = 0.4107214552505143
$ws.Range("M5").Value = 8.841467
$ws.Range("N5").Value = 26.524401
$ws.Range("O5").Value = 0.5917001192060068
$ws.Range("P5").Value = 0.5917001192060067
$ws.Range("Q5").Value = 60.03247637671466
$ws.Range("R5").Value = 540.292287390432
$ws.Range("S5").Value = 0.243023934032194
$ws.Range("T5").Value = 0.2430239340321939

# Row 6
$ws.Range("I6").Value = 0.4107214552505144
$ws.Range("J6").Value = 0.4107214552505143
$ws.Range("O6").Value = 0.2746155987184545
$ws.Range("P6").Value = 0.2746155987184545
$ws.Range("S6").Value = 0.1127905183401349
$ws.Range("T6").Value = 0.1127905183401349

# Row 7
$ws.Range("I7").Value = 0.4107214552505144
$ws.Range("J7").Value = 0.4107214552505143
$ws.Range("M7").Value = 1.997574666666667
$ws.Range("N7").Value = 5.992724
$ws.Range("O7").Value = 0.1336842820755386
$ws.Range("P7").Value = 0.1336842820755386
$ws.Range("Q7").Value = 13.56328695084089
$ws.Range("R7").Value = 122.069582557568
$ws.Range("S7").Value = 0.05490700287818548
$ws.Range("T7").Value = 0.05490700287818545

# Row 8
$ws.Range("G8").Value = 3.710753333333333
$ws.Range("H8").Value = 11.13226
$ws.Range("I8").Value = 0.2244644393883547
$ws.Range("J8").Value = 0.2244644393883547
$ws.Range("M8").Value = 8.841467
$ws.Range("N8").Value = 26.524401
$ws.Range("O8").Value = 0.5917001192060068
$ws.Range("P8").Value = 0.5917001192060067
$ws.Range("Q8").Value = 32.80850314180666
$ws.Range("R8").Value = 295.27652827626
$ws.Range("S8").Value = 0.132815635543599
$ws.Range("T8").Value = 0.1328156355435989

# Row 9
$ws.Range("G9").Value = 3.710753333333333
$ws.Range("H9").Value = 11.13226
$ws.Range("I9").Value = 0.2244644393883547
$ws.Range("J9").Value = 0.2244644393883547
$ws.Range("O9").Value = 0.2746155987184545
$ws.Range("P9").Value = 0.2746155987184545
$ws.Range("Q9").Value = 15.22684623662666
$ws.Range("R9").Value = 137.04161612964
$ws.Range("S9").Value = 0.06164143641363528
$ws.Range("T9").Value = 0.06164143641363525

# Row 10
$ws.Range("G10").Value = 3.710753333333333
$ws.Range("H10").Value = 11.13226
$ws.Range("I10").Value = 0.2244644393883547
$ws.Range("J10").Value = 0.2244644393883547
$ws.Range("M10").Value = 1.997574666666667
$ws.Range("N10").Value = 5.992724
$ws.Range("O10").Value = 0.1336842820755386
$ws.Range("P10").Value = 0.1336842820755386
$ws.Range("Q10").Value = 7.412506852915555
$ws.Range("R10").Value = 66.71256167623999
$ws.Range("S10").Value = 0.03000736743112046
$ws.Range("T10").Value = 0.03000736743112044
